$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (English item keys) - new columns AM..BL
$ws.Cells.Item(1, 39).Value = "fast_shoes"
$ws.Cells.Item(1, 40).Value = "stab_shield"
$ws.Cells.Item(1, 41).Value = "arm_shield"
$ws.Cells.Item(1, 42).Value = "wood_shield"
$ws.Cells.Item(1, 43).Value = "long_arch"
$ws.Cells.Item(1, 44).Value = "wood_arch"
$ws.Cells.Item(1, 45).Value = "zhuge_bow"
$ws.Cells.Item(1, 46).Value = "multi_bow"
$ws.Cells.Item(1, 47).Value = "hand_bow"
$ws.Cells.Item(1, 48).Value = "shadow_charm"
$ws.Cells.Item(1, 49).Value = "exchange_charm"
$ws.Cells.Item(1, 50).Value = "silent_charm"
$ws.Cells.Item(1, 51).Value = "confine_charm"
$ws.Cells.Item(1, 52).Value = "thunder_charm"
$ws.Cells.Item(1, 53).Value = "king_arrow"
$ws.Cells.Item(1, 54).Value = "treasure_bowl"
$ws.Cells.Item(1, 55).Value = "shield_token"
$ws.Cells.Item(1, 56).Value = "sword_stone"
$ws.Cells.Item(1, 57).Value = "energy_potion"
$ws.Cells.Item(1, 58).Value = "angry_potion"
$ws.Cells.Item(1, 59).Value = "lucky_potion"
$ws.Cells.Item(1, 60).Value = "boutique_feather"
$ws.Cells.Item(1, 61).Value = "boutique_rice"
$ws.Cells.Item(1, 62).Value = "trophy"
$ws.Cells.Item(1, 63).Value = "jerky"
$ws.Cells.Item(1, 64).Value = "month_card"

# Row 2 (Chinese item names) - new columns AM..BL
$ws.Cells.Item(2, 39).Value = "极速靴"
$ws.Cells.Item(2, 40).Value = "刺盾"
$ws.Cells.Item(2, 41).Value = "军盾"
$ws.Cells.Item(2, 42).Value = "木盾"
$ws.Cells.Item(2, 43).Value = "长弓"
$ws.Cells.Item(2, 44).Value = "木弓"
$ws.Cells.Item(2, 45).Value = "诸葛连弩"
$ws.Cells.Item(2, 46).Value = "连射弩"
$ws.Cells.Item(2, 47).Value = "手弩"
$ws.Cells.Item(2, 48).Value = "幻影符"
$ws.Cells.Item(2, 49).Value = "交换符"
$ws.Cells.Item(2, 50).Value = "静默符"
$ws.Cells.Item(2, 51).Value = "沉默符"
$ws.Cells.Item(2, 52).Value = "奔雷符"
$ws.Cells.Item(2, 53).Value = "王之箭矢"
$ws.Cells.Item(2, 54).Value = "聚宝盆"
$ws.Cells.Item(2, 55).Value = "护盾令"
$ws.Cells.Item(2, 56).Value = "剑意石"
$ws.Cells.Item(2, 57).Value = "精力药水"
$ws.Cells.Item(2, 58).Value = "愤怒药水"
$ws.Cells.Item(2, 59).Value = "幸运药水"
$ws.Cells.Item(2, 60).Value = "精致彩羽"
$ws.Cells.Item(2, 61).Value = "精品稻米"
$ws.Cells.Item(2, 62).Value = "奖杯"
$ws.Cells.Item(2, 63).Value = "风干肉条"
$ws.Cells.Item(2, 64).Value = "月卡"

# Update selection and scroll position to match final state
$ws.Range("K8").Select()
try {
    $win = $excel.ActiveWindow
    $win.ScrollColumn = 4
    $win.ScrollRow = 1
} catch {
    # scrolling the viewport is a cosmetic, best-effort operation only
}
